$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text formatting
# (values like "93.70" or "0.00001056" would otherwise be auto-converted to numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.649.50"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.775.46"
$ws.Range("E3").Value = "  +1.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.74"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4622"
$ws.Range("E7").Value = "  +4.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3581"
$ws.Range("E8").Value = "  -0.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07475"
$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.74"
$ws.Range("E10").Value = "  -1.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.099"
$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9989"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.86"
$ws.Range("E13").Value = "  +1.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.045"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.206"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.772.97"
$ws.Range("E16").Value = "  +1.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.70"
$ws.Range("E17").Value = "  +1.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001056"
$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06425"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9985"
$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.10"
$ws.Range("E21").Value = "  +1.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.726.46"
$ws.Range("E23").Value = "  +0.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.29"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.082"
$ws.Range("E25").Value = "  -0.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.85"
$ws.Range("E26").Value = "  +2.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.26"
$ws.Range("E27").Value = "  -0.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.974.34"
$ws.Range("E28").Value = "  +1.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.160"
$ws.Range("E29").Value = "  +2.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.09"
$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.102"
$ws.Range("E31").Value = "  +2.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09208"
$ws.Range("E32").Value = "  +2.86%  "

$ws.Range("E33").Value = "  +0.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.544"
$ws.Range("E34").Value = "  +0.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.82"
$ws.Range("E35").Value = "  -1.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02289"
$ws.Range("E36").Value = "  -0.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06110"
$ws.Range("E37").Value = "  +2.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2092"
$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6306"
$ws.Range("E39").Value = "  -0.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.956"
$ws.Range("E40").Value = "  +0.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.183"
$ws.Range("E41").Value = "  -1.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.391"
$ws.Range("E42").Value = "  +0.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.804"
$ws.Range("E43").Value = "  +0.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.20"
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("E45").Value = "  +0.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5887"
$ws.Range("E46").Value = "  +0.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.38"
$ws.Range("E47").Value = "  +0.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.949"
$ws.Range("E48").Value = "  +0.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06943"
$ws.Range("E49").Value = "  +1.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.136"
$ws.Range("E50").Value = "  -1.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.46"
$ws.Range("E51").Value = "  +0.69%  "
